# Refresh the hourly cryptos snapshot (Price/Volume columns) to match the
# latest scrape, per the "Updated cryptos list ... GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "25.253.60"
# Row 3: Ethereum
$ws.Range("D3").Value = "1.554.22"
# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.13%  "
# Row 5: BNB
$ws.Range("D5").Value = "'207.08"
$ws.Range("E5").Value = "  -3.32%  "
# Row 6: USDC
$ws.Range("E6").Value = "  -0.14%  "
# Row 7: XRP
$ws.Range("D7").Value = "'0.476"
$ws.Range("E7").Value = "  -5.40%  "
# Row 8: Dogecoin
$ws.Range("D8").Value = "'0.0608"
$ws.Range("E8").Value = "  -1.75%  "
# Row 9: Cardano
$ws.Range("E9").Value = "  -3.35%  "
# Row 10: Solana
$ws.Range("D10").Value = "'17.69"
$ws.Range("E10").Value = "  -4.27%  "
# Row 11: TRON
$ws.Range("D11").Value = "'0.0780"
$ws.Range("E11").Value = "  -1.06%  "
# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.770.32"
# Row 13: WrappedEther
$ws.Range("D13").Value = "1.550.25"
$ws.Range("E13").Value = "  -4.94%  "
# Row 14: Polkadot
$ws.Range("E14").Value = "  -4.40%  "
# Row 15: Polygon
$ws.Range("D15").Value = "'0.504"
$ws.Range("E15").Value = "  -4.36%  "
# Row 16: WrappedBTC
$ws.Range("D16").Value = "25.260.73"
$ws.Range("E16").Value = "  -2.94%  "
# Row 17: Litecoin
$ws.Range("D17").Value = "'58.77"
# Row 18: ShibaInu
$ws.Range("E18").Value = "  -4.76%  "
# Row 19: Dai
$ws.Range("E19").Value = "  -0.11%  "
# Row 20: BitcoinCash
$ws.Range("D20").Value = "'185.39"
$ws.Range("E20").Value = "  -3.76%  "
# Row 21: Uniswap
$ws.Range("E21").Value = "  -3.48%  "
# Row 22: Avalanche
$ws.Range("E22").Value = "  -2.82%  "
# Row 23: Chainlink
$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  -3.71%  "
# Row 24: Stellar
$ws.Range("E24").Value = "  -4.10%  "
# Row 25: BinanceUSD
$ws.Range("E25").Value = "  -0.08%  "
# Row 26: Monero
$ws.Range("D26").Value = "'140.43"
$ws.Range("E26").Value = "  -2.73%  "
# Row 27: Toncoin
$ws.Range("E27").Value = "  -4.92%  "
# Row 28: EthereumClassic
$ws.Range("D28").Value = "'14.84"
# Row 29: Cosmos
$ws.Range("E29").Value = "  -4.95%  "
# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -6.86%  "
# Row 31: Hedera
$ws.Range("D31").Value = "'0.0466"
$ws.Range("E31").Value = "  -3.39%  "
# Row 32: Filecoin
$ws.Range("D32").Value = "'3.02"
$ws.Range("E32").Value = "  -3.29%  "
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "  -5.05%  "
# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  -3.07%  "
# Row 35: HuobiToken
$ws.Range("E35").Value = "  -3.87%  "
# Row 36: Maker
$ws.Range("D36").Value = "1.083.08"
$ws.Range("E36").Value = "  -4.00%  "
# Row 37: PaxDollar
$ws.Range("E37").Value = "  -0.15%  "
# Row 38: VeChain
$ws.Range("E38").Value = "  -3.14%  "
# Row 39: ImmutableX
$ws.Range("E39").Value = "  -4.65%  "
# Row 40: ARBITRUM
$ws.Range("E40").Value = "  -10.15%  "
# Row 41: MXToken
$ws.Range("E41").Value = "  -7.52%  "
# Row 42: TrustWalletToken
$ws.Range("D42").Value = "'0.797"
$ws.Range("E42").Value = "  +5.47%  "
# Row 43: Quant
$ws.Range("D43").Value = "'92.46"
$ws.Range("E43").Value = "  -5.92%  "
# Row 44: FraxShare
$ws.Range("D44").Value = "'5.04"
$ws.Range("E44").Value = "  -1.55%  "
# Row 45: RocketPoolETH
$ws.Range("D45").Value = "1.685.45"
$ws.Range("E45").Value = "  -4.57%  "
# Row 46: BabyDogeCoin
$ws.Range("E46").Value = "  -3.05%  "
# Row 47: RenderToken
$ws.Range("E47").Value = "  -1.84%  "
# Row 48: Aave
$ws.Range("E48").Value = "  -3.87%  "
# Row 49: Cronos
$ws.Range("D49").Value = "'0.0501"
$ws.Range("E49").Value = "  -4.24%  "
# Row 50: USDD
$ws.Range("E50").Value = "  -0.30%  "
# Row 51: Mantle
$ws.Range("E51").Value = "  -2.06%  "
